$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 3 ("mid group" / "N"), pushing the existing
# rows 3-8 ("last group" ... "marksheet") down to 4-9.
$ws.Rows("3").Insert()
$ws.Range("A3").Value = "mid group"
$ws.Range("B3").Value = "N"
$ws.Rows("3").RowHeight = 19

# The "questions per category" value cell (now at B7) becomes
# left-aligned in the new layout.
$ws.Range("B7").HorizontalAlignment = -4131

# Move the active selection to B2.
$ws.Range("B2").Select() | Out-Null
